# Add team record columns (Wins, Losses, Ties) to the NYM_2009 data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (AC1) onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record for every data row (rows 2 through 55).
$ws.Range("AD2:AD55").Value = 70
$ws.Range("AE2:AE55").Value = 92
$ws.Range("AF2:AF55").Value = 0
